# Weekly update: insert this week's "Betarraga" price records at the top
# of the data block (row 501), pushing all older rows down by 3 rows.
#
# The sheet holds one row per (Fecha, Calidad) combination, newest entries
# first within the date-ordered block that starts at row 501. Inserting 3
# full rows there (one each for "Primera", "Segunda", "Tercera") shifts the
# remaining historical rows from 501:586 down to 504:589, growing the used
# range from A1:R586 to A1:R589 - matching the diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("501:503").Insert()

# New week's data: Fecha serial 44522, three quality grades.
$newRows = @(
    @{ Row=501; Calidad="Primera"; Volumen=51000; PMin=85; PMax=90; PProm=87 },
    @{ Row=502; Calidad="Segunda"; Volumen=42000; PMin=65; PMax=70; PProm=67 },
    @{ Row=503; Calidad="Tercera"; Volumen=14000; PMin=40; PMax=40; PProm=40 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 6
    $ws.Cells.Item($row, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value  = "Metropolitana"
    $ws.Cells.Item($row, 4).Value  = 44522
    $ws.Cells.Item($row, 5).Value  = 13
    $ws.Cells.Item($row, 6).Value  = 100114014
    $ws.Cells.Item($row, 7).Value  = "Betarraga"
    $ws.Cells.Item($row, 8).Value  = "Sin especificar"
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = "$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
